$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Fitness") values for rows 2..202 (Generation 0..200) become 7534
$ws.Range("C2:C202").Value = 7534

# Column C ("Fitness") values for rows 203..252 (Generation 201..250) become 7295
$ws.Range("C203:C252").Value = 7295
